$wb = $excel.ActiveWorkbook

# --- Rename header cells -------------------------------------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the existing sheets ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy over formatting (header style + date style) from an existing sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A27").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
    @(44983.99999999999, 3, -7.366073852294731, 11.70402306038798),
    @(45004.99999999999, 3, -7.216849832438281, 11.82536587524062),
    @(45032.99999999999, 3, -7.582955310303777, 13.51294457330778),
    @(45039.99999999999, 3, -6.021279566024782, 13.27920897604844),
    @(45340.99999999999, 6, -3.885058266441594, 15.93474944122173),
    @(45354.99999999999, 6, -3.792693546850624, 15.66652363172814),
    @(45368.99999999999, 6, -4.57610600148247,  15.81503029297246),
    @(45382.99999999999, 6, -3.530528320994505, 16.62817076881087),
    @(45396.99999999999, 6, -3.557612168012401, 16.96853349438583),
    @(45459.99999999999, 7, -3.375406494570448, 16.67206084318449),
    @(45487.99999999999, 7, -2.547092556123525, 17.17344313702716),
    @(45508.99999999999, 7, -3.055858215753123, 17.06581873242048),
    @(45515.99999999999, 7, -2.530142862552939, 16.97925756431568),
    @(45550.99999999999, 8, -2.117899754495333, 17.26704739381196),
    @(45564.99999999999, 8, -1.67983065440064,  16.47469043119569),
    @(45578.99999999999, 8, -2.189458617485215, 17.937428389936),
    @(45585.99999999999, 8, -1.786539597259435, 17.77741193339039),
    @(45592.99999999999, 8, -1.939381168272476, 18.32008119852734),
    @(45599.99999999999, 8, -2.346072482106258, 18.50177633451825),
    @(45606.99999999999, 8, -1.285092315723528, 18.4266726859014),
    @(45613.99999999999, 8, -1.262469193470042, 17.86282720944297),
    @(45620.99999999999, 8, -1.569581324880331, 17.91372017809743),
    @(45627.99999999999, 8, -1.668182183836023, 18.12673230858715),
    @(45634.99999999999, 8, -1.48392628707545,  18.14702367429435),
    @(45641.99999999999, 8, -1.382246057125505, 18.47173708257656),
    @(45648.99999999999, 9, -1.330328867163401, 18.28700458589067)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Restore original active sheet/selection so view-state stays unchanged
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null
